$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 111469967
$ws.Range("B10").Value = 5113
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 100526
$ws.Range("F10").Value = 'Bronshjon'
$ws.Range("G10").Value = 'Callidium coriaceum'
$ws.Range("H10").Value = 'Paykull, 1800'
$ws.Range("J10").ClearContents() | Out-Null
$ws.Range("K10").ClearContents() | Out-Null
$ws.Range("L10").ClearContents() | Out-Null
$ws.Range("M10").Value = 'äldre gnagspår'
$ws.Range("N10").ClearContents() | Out-Null
$ws.Range("Q10").Value = 554705.6319759471
$ws.Range("R10").Value = 6698113.601669285
$ws.Range("AF10").ClearContents() | Out-Null

# Row 11
$ws.Range("A11").Value = 111469965
$ws.Range("B11").Value = 5113
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 100526
$ws.Range("F11").Value = 'Bronshjon'
$ws.Range("G11").Value = 'Callidium coriaceum'
$ws.Range("H11").Value = 'Paykull, 1800'
$ws.Range("J11").ClearContents() | Out-Null
$ws.Range("K11").ClearContents() | Out-Null
$ws.Range("L11").ClearContents() | Out-Null
$ws.Range("M11").Value = 'färska gnagspår'
$ws.Range("N11").ClearContents() | Out-Null
$ws.Range("Q11").Value = 554716.1509068209
$ws.Range("R11").Value = 6698137.967376946
$ws.Range("AF11").ClearContents() | Out-Null

# Row 12
$ws.Range("A12").Value = 111469966
$ws.Range("B12").Value = 5113
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 100526
$ws.Range("F12").Value = 'Bronshjon'
$ws.Range("G12").Value = 'Callidium coriaceum'
$ws.Range("H12").Value = 'Paykull, 1800'
$ws.Range("J12").ClearContents() | Out-Null
$ws.Range("K12").ClearContents() | Out-Null
$ws.Range("L12").ClearContents() | Out-Null
$ws.Range("M12").Value = 'äldre gnagspår'
$ws.Range("N12").ClearContents() | Out-Null
$ws.Range("Q12").Value = 554729.2459973614
$ws.Range("R12").Value = 6698057.144588907
$ws.Range("AF12").ClearContents() | Out-Null

# Row 13
$ws.Range("A13").Value = 111469944
$ws.Range("B13").Value = 96348
$ws.Range("D13").Value = 'VU'
$ws.Range("E13").Value = 220787
$ws.Range("F13").Value = 'Knärot'
$ws.Range("G13").Value = 'Goodyera repens'
$ws.Range("H13").Value = '(L.) R. Br.'
$ws.Range("J13").ClearContents() | Out-Null
$ws.Range("K13").ClearContents() | Out-Null
$ws.Range("L13").ClearContents() | Out-Null
$ws.Range("M13").ClearContents() | Out-Null
$ws.Range("N13").ClearContents() | Out-Null
$ws.Range("Q13").Value = 554647.0313377964
$ws.Range("R13").Value = 6697996.127116338
$ws.Range("AF13").ClearContents() | Out-Null

# Row 14
$ws.Range("A14").Value = 111469926
$ws.Range("B14").Value = 89369
$ws.Range("D14").Value = 'LC'
$ws.Range("E14").Value = 5447
$ws.Range("F14").Value = 'Vedticka'
$ws.Range("G14").Value = 'Fuscoporia viticola'
$ws.Range("H14").Value = '(Schwein.) Murrill'
$ws.Range("J14").ClearContents() | Out-Null
$ws.Range("K14").ClearContents() | Out-Null
$ws.Range("L14").ClearContents() | Out-Null
$ws.Range("M14").ClearContents() | Out-Null
$ws.Range("N14").ClearContents() | Out-Null
$ws.Range("Q14").Value = 554745.7538377594
$ws.Range("R14").Value = 6698078.142900295
$ws.Range("AF14").ClearContents() | Out-Null

# Row 15
$ws.Range("A15").Value = 111469950
$ws.Range("B15").Value = 96348
$ws.Range("D15").Value = 'VU'
$ws.Range("E15").Value = 220787
$ws.Range("F15").Value = 'Knärot'
$ws.Range("G15").Value = 'Goodyera repens'
$ws.Range("H15").Value = '(L.) R. Br.'
$ws.Range("J15").ClearContents() | Out-Null
$ws.Range("K15").ClearContents() | Out-Null
$ws.Range("L15").ClearContents() | Out-Null
$ws.Range("M15").ClearContents() | Out-Null
$ws.Range("N15").ClearContents() | Out-Null
$ws.Range("Q15").Value = 554648.2514272946
$ws.Range("R15").Value = 6697980.830233379
$ws.Range("AF15").ClearContents() | Out-Null

# Row 16
$ws.Range("A16").Value = 111469964
$ws.Range("B16").Value = 5113
$ws.Range("D16").Value = 'LC'
$ws.Range("E16").Value = 100526
$ws.Range("F16").Value = 'Bronshjon'
$ws.Range("G16").Value = 'Callidium coriaceum'
$ws.Range("H16").Value = 'Paykull, 1800'
$ws.Range("J16").ClearContents() | Out-Null
$ws.Range("K16").ClearContents() | Out-Null
$ws.Range("L16").ClearContents() | Out-Null
$ws.Range("M16").Value = 'färska gnagspår'
$ws.Range("N16").ClearContents() | Out-Null
$ws.Range("Q16").Value = 554731.9372321201
$ws.Range("R16").Value = 6698141.169601779
$ws.Range("AF16").ClearContents() | Out-Null

# Row 17
$ws.Range("A17").Value = 111469962
$ws.Range("B17").Value = 5113
$ws.Range("D17").Value = 'LC'
$ws.Range("E17").Value = 100526
$ws.Range("F17").Value = 'Bronshjon'
$ws.Range("G17").Value = 'Callidium coriaceum'
$ws.Range("H17").Value = 'Paykull, 1800'
$ws.Range("J17").ClearContents() | Out-Null
$ws.Range("K17").ClearContents() | Out-Null
$ws.Range("L17").ClearContents() | Out-Null
$ws.Range("M17").Value = 'färska gnagspår'
$ws.Range("N17").ClearContents() | Out-Null
$ws.Range("Q17").Value = 554640.2091243146
$ws.Range("R17").Value = 6697989.107814683
$ws.Range("AF17").ClearContents() | Out-Null

# Row 18
$ws.Range("A18").Value = 111469941
$ws.Range("B18").Value = 96348
$ws.Range("D18").Value = 'VU'
$ws.Range("E18").Value = 220787
$ws.Range("F18").Value = 'Knärot'
$ws.Range("G18").Value = 'Goodyera repens'
$ws.Range("H18").Value = '(L.) R. Br.'
$ws.Range("J18").ClearContents() | Out-Null
$ws.Range("K18").Value = 'blomning'
$ws.Range("L18").ClearContents() | Out-Null
$ws.Range("M18").ClearContents() | Out-Null
$ws.Range("N18").ClearContents() | Out-Null
$ws.Range("Q18").Value = 554704.8063610581
$ws.Range("R18").Value = 6698102.720679003
$ws.Range("AF18").ClearContents() | Out-Null

# Row 19
$ws.Range("A19").Value = 111469953
$ws.Range("B19").Value = 96348
$ws.Range("D19").Value = 'VU'
$ws.Range("E19").Value = 220787
$ws.Range("F19").Value = 'Knärot'
$ws.Range("G19").Value = 'Goodyera repens'
$ws.Range("H19").Value = '(L.) R. Br.'
$ws.Range("J19").ClearContents() | Out-Null
$ws.Range("K19").ClearContents() | Out-Null
$ws.Range("L19").ClearContents() | Out-Null
$ws.Range("M19").ClearContents() | Out-Null
$ws.Range("N19").ClearContents() | Out-Null
$ws.Range("Q19").Value = 554668.8331894471
$ws.Range("R19").Value = 6698027.085862564
$ws.Range("AF19").ClearContents() | Out-Null

# Row 20
$ws.Range("A20").Value = 111469952
$ws.Range("B20").Value = 96348
$ws.Range("D20").Value = 'VU'
$ws.Range("E20").Value = 220787
$ws.Range("F20").Value = 'Knärot'
$ws.Range("G20").Value = 'Goodyera repens'
$ws.Range("H20").Value = '(L.) R. Br.'
$ws.Range("J20").ClearContents() | Out-Null
$ws.Range("K20").ClearContents() | Out-Null
$ws.Range("L20").ClearContents() | Out-Null
$ws.Range("M20").ClearContents() | Out-Null
$ws.Range("N20").ClearContents() | Out-Null
$ws.Range("Q20").Value = 554701.1291447466
$ws.Range("R20").Value = 6697985.57934437
$ws.Range("AF20").ClearContents() | Out-Null

# Row 21
$ws.Range("A21").Value = 111469963
$ws.Range("B21").Value = 5113
$ws.Range("D21").Value = 'LC'
$ws.Range("E21").Value = 100526
$ws.Range("F21").Value = 'Bronshjon'
$ws.Range("G21").Value = 'Callidium coriaceum'
$ws.Range("H21").Value = 'Paykull, 1800'
$ws.Range("J21").ClearContents() | Out-Null
$ws.Range("K21").ClearContents() | Out-Null
$ws.Range("L21").ClearContents() | Out-Null
$ws.Range("M21").Value = 'färska gnagspår'
$ws.Range("N21").ClearContents() | Out-Null
$ws.Range("Q21").Value = 554718.6790950731
$ws.Range("R21").Value = 6698003.135367867
$ws.Range("AF21").ClearContents() | Out-Null

# Row 22
$ws.Range("A22").Value = 111469922
$ws.Range("B22").Value = 5135
$ws.Range("D22").Value = 'LC'
$ws.Range("E22").Value = 105930
$ws.Range("F22").Value = 'Vågbandad barkbock'
$ws.Range("G22").Value = 'Semanotus undatus'
$ws.Range("H22").Value = '(Linnaeus, 1758)'
$ws.Range("J22").ClearContents() | Out-Null
$ws.Range("K22").ClearContents() | Out-Null
$ws.Range("L22").ClearContents() | Out-Null
$ws.Range("M22").Value = 'äldre gnagspår'
$ws.Range("N22").ClearContents() | Out-Null
$ws.Range("Q22").Value = 554716.6256586342
$ws.Range("R22").Value = 6698008.044787553
$ws.Range("AF22").ClearContents() | Out-Null

# Row 23
$ws.Range("A23").Value = 111469968
$ws.Range("B23").Value = 5113
$ws.Range("D23").Value = 'LC'
$ws.Range("E23").Value = 100526
$ws.Range("F23").Value = 'Bronshjon'
$ws.Range("G23").Value = 'Callidium coriaceum'
$ws.Range("H23").Value = 'Paykull, 1800'
$ws.Range("J23").ClearContents() | Out-Null
$ws.Range("K23").ClearContents() | Out-Null
$ws.Range("L23").ClearContents() | Out-Null
$ws.Range("M23").Value = 'äldre gnagspår'
$ws.Range("N23").ClearContents() | Out-Null
$ws.Range("Q23").Value = 554679.218646974
$ws.Range("R23").Value = 6698060.342582431
$ws.Range("AF23").ClearContents() | Out-Null

# Row 24
$ws.Range("A24").Value = 111469946
$ws.Range("B24").Value = 96348
$ws.Range("D24").Value = 'VU'
$ws.Range("E24").Value = 220787
$ws.Range("F24").Value = 'Knärot'
$ws.Range("G24").Value = 'Goodyera repens'
$ws.Range("H24").Value = '(L.) R. Br.'
$ws.Range("J24").ClearContents() | Out-Null
$ws.Range("K24").ClearContents() | Out-Null
$ws.Range("L24").ClearContents() | Out-Null
$ws.Range("M24").ClearContents() | Out-Null
$ws.Range("N24").ClearContents() | Out-Null
$ws.Range("Q24").Value = 554664.6782300239
$ws.Range("R24").Value = 6698007.261790544
$ws.Range("AF24").ClearContents() | Out-Null

# Row 25
$ws.Range("A25").Value = 111469958
$ws.Range("B25").Value = 89621
$ws.Range("D25").Value = 'NT'
$ws.Range("E25").Value = 1101
$ws.Range("F25").Value = 'Gropticka'
$ws.Range("G25").Value = 'Postia guttulata'
$ws.Range("H25").Value = '(Peck) Jülich'
$ws.Range("J25").ClearContents() | Out-Null
$ws.Range("K25").ClearContents() | Out-Null
$ws.Range("L25").ClearContents() | Out-Null
$ws.Range("M25").ClearContents() | Out-Null
$ws.Range("N25").ClearContents() | Out-Null
$ws.Range("Q25").Value = 554681.1975678616
$ws.Range("R25").Value = 6698060.372405332
$ws.Range("AF25").ClearContents() | Out-Null

# Row 26
$ws.Range("A26").Value = 111469951
$ws.Range("B26").Value = 96348
$ws.Range("D26").Value = 'VU'
$ws.Range("E26").Value = 220787
$ws.Range("F26").Value = 'Knärot'
$ws.Range("G26").Value = 'Goodyera repens'
$ws.Range("H26").Value = '(L.) R. Br.'
$ws.Range("J26").ClearContents() | Out-Null
$ws.Range("K26").ClearContents() | Out-Null
$ws.Range("L26").ClearContents() | Out-Null
$ws.Range("M26").ClearContents() | Out-Null
$ws.Range("N26").ClearContents() | Out-Null
$ws.Range("Q26").Value = 554679.0891228422
$ws.Range("R26").Value = 6697970.425878088
$ws.Range("AF26").ClearContents() | Out-Null

# Row 27
$ws.Range("A27").Value = 111469954
$ws.Range("B27").Value = 96348
$ws.Range("D27").Value = 'VU'
$ws.Range("E27").Value = 220787
$ws.Range("F27").Value = 'Knärot'
$ws.Range("G27").Value = 'Goodyera repens'
$ws.Range("H27").Value = '(L.) R. Br.'
$ws.Range("J27").ClearContents() | Out-Null
$ws.Range("K27").ClearContents() | Out-Null
$ws.Range("L27").ClearContents() | Out-Null
$ws.Range("M27").ClearContents() | Out-Null
$ws.Range("N27").ClearContents() | Out-Null
$ws.Range("Q27").Value = 554709.4759112563
$ws.Range("R27").Value = 6698022.75809369
$ws.Range("AF27").ClearContents() | Out-Null

# Row 28
$ws.Range("A28").Value = 111469949
$ws.Range("B28").Value = 96348
$ws.Range("D28").Value = 'VU'
$ws.Range("E28").Value = 220787
$ws.Range("F28").Value = 'Knärot'
$ws.Range("G28").Value = 'Goodyera repens'
$ws.Range("H28").Value = '(L.) R. Br.'
$ws.Range("J28").ClearContents() | Out-Null
$ws.Range("K28").ClearContents() | Out-Null
$ws.Range("L28").ClearContents() | Out-Null
$ws.Range("M28").ClearContents() | Out-Null
$ws.Range("N28").ClearContents() | Out-Null
$ws.Range("Q28").Value = 554654.1362404823
$ws.Range("R28").Value = 6697984.37715952
$ws.Range("AF28").ClearContents() | Out-Null

# Row 29
$ws.Range("A29").Value = 111469947
$ws.Range("B29").Value = 96348
$ws.Range("D29").Value = 'VU'
$ws.Range("E29").Value = 220787
$ws.Range("F29").Value = 'Knärot'
$ws.Range("G29").Value = 'Goodyera repens'
$ws.Range("H29").Value = '(L.) R. Br.'
$ws.Range("J29").ClearContents() | Out-Null
$ws.Range("K29").ClearContents() | Out-Null
$ws.Range("L29").ClearContents() | Out-Null
$ws.Range("M29").ClearContents() | Out-Null
$ws.Range("N29").ClearContents() | Out-Null
$ws.Range("Q29").Value = 554660.8096201464
$ws.Range("R29").Value = 6698001.275046931
$ws.Range("AF29").ClearContents() | Out-Null

# Row 30
$ws.Range("A30").Value = 111469929
$ws.Range("B30").Value = 56398
$ws.Range("D30").Value = 'NT'
$ws.Range("E30").Value = 100109
$ws.Range("F30").Value = 'Tretåig hackspett'
$ws.Range("G30").Value = 'Picoides tridactylus'
$ws.Range("H30").Value = '(Linnaeus, 1758)'
$ws.Range("J30").ClearContents() | Out-Null
$ws.Range("K30").ClearContents() | Out-Null
$ws.Range("L30").ClearContents() | Out-Null
$ws.Range("M30").Value = 'äldre spår'
$ws.Range("N30").ClearContents() | Out-Null
$ws.Range("Q30").Value = 554646.3468513897
$ws.Range("R30").Value = 6697975.861129273
$ws.Range("AF30").ClearContents() | Out-Null
